$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("packet layout")

# --- Shared strings: remove "OPTIONAL" marker from rows 21-28 (column D) ---
$ws.Range("D21:D28").ClearContents()

# --- Re-point B column labels (EFM/Battery/PD/Cloud mean/etc. unaffected in meaning,
#     they're simply re-indexed in the underlying string table once OPTIONAL goes away) ---
# (no textual change needed for B column -- values already correct)

# --- Add new Status Byte bitfield table (H1:P2) ---
$ws.Range("H1").Value = "Status Byte"
$ws.Range("I1").Value = 7
$ws.Range("J1").Value = 6
$ws.Range("K1").Value = 5
$ws.Range("L1").Value = 4
$ws.Range("M1").Value = 3
$ws.Range("N1").Value = 2
$ws.Range("O1").Value = 1
$ws.Range("P1").Value = 0
$ws.Range("I1:P1").Font.Bold = $true

$ws.Range("H2").Value = "Meaning"
$ws.Range("H2").Font.Bold = $true
$ws.Range("P2").Value = "Promiscuous mode"
$ws.Range("O2").Value = "GPS valid"
$ws.Range("N2").Value = "GPS >3 sat"
$ws.Range("M2").Value = "GPS hdop ok"
$ws.Range("L2").Value = "SD init ok"

$ws.Columns.Item(16).AutoFit() | Out-Null

# --- Selection / view changes ---
$ws.Range("D20:D28").Select()

# --- Workbook window height tweak ---
$excel.ActiveWindow.Height = 14300
